# Deploying to gh-pages: add 2021 column (O) to the 9.3.1 small-scale
# industries table and update the 2018/2020 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend formatting of the last existing year-column (N) into the new
#    2021 column (O) for the header divider row, the year-label row and the
#    data row, so the new column inherits borders/fonts/alignment.
$ws.Range("N3:N5").Copy() | Out-Null
$ws.Range("O3:O5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Fill in the new 2021 column's values.
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 4.0999999999999996

# 3) Corrected historical figures.
$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1

# 4) Move the active selection the same way the original author's session
#    ended up (row 4 instead of row 6) after the edit.
$ws.Range("P4").Select() | Out-Null
